# Adds three new worksheets to the workbook, mirroring the structure that a
# companion "statistics" Cypher query run produces alongside the existing
# CypherOutput / Message sheets:
#   - CypherOutput_Message : exact copy of the "Message" sheet
#   - StatOutput            : header row + one data row with the stat counts
#   - StatOutput_Message    : like "Message" but the embedded Cypher text is
#                              the stats query instead of the original one

$wb = $excel.ActiveWorkbook
$msgSheet = $wb.Worksheets.Item("Message")

$statQuery = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN[''Stage 2'']  OPTIONAL MATCH (f:file)-[*]-->(c), (samp:sample)-[*]-->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'

# ---- 1) CypherOutput_Message : straight copy of Message ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cypherMsgSheet = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$cypherMsgSheet.Name = "CypherOutput_Message"
for ($r = 1; $r -le 10; $r++) {
    $cypherMsgSheet.Cells.Item($r, 1).Value = $msgSheet.Cells.Item($r, 1).Value2
}

# ---- 2) StatOutput : headers + counts row ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$statSheet = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$statSheet.Name = "StatOutput"
$statSheet.Cells.Item(1, 1).Value = "number_of_files"
$statSheet.Cells.Item(1, 2).Value = "number_of_sample"
$statSheet.Cells.Item(1, 3).Value = "number_of_cases"
$statSheet.Cells.Item(1, 4).Value = "number_of_study"
# Values are written as text (quote-prefixed) to match the source data,
# which stores these counts as text rather than numbers.
$statSheet.Cells.Item(2, 1).Value = "'0"
$statSheet.Cells.Item(2, 2).Value = "'0"
$statSheet.Cells.Item(2, 3).Value = "'1"
$statSheet.Cells.Item(2, 4).Value = "'1"

# ---- 3) StatOutput_Message : Message content twice, with the 8th line of
#         the second copy replaced by the stats Cypher query ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$statMsgSheet = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$statMsgSheet.Name = "StatOutput_Message"
for ($r = 1; $r -le 10; $r++) {
    $statMsgSheet.Cells.Item($r, 1).Value = $msgSheet.Cells.Item($r, 1).Value2
}
for ($r = 1; $r -le 7; $r++) {
    $statMsgSheet.Cells.Item($r + 10, 1).Value = $msgSheet.Cells.Item($r, 1).Value2
}
$statMsgSheet.Cells.Item(18, 1).Value = $statQuery
for ($r = 9; $r -le 10; $r++) {
    $statMsgSheet.Cells.Item($r + 10, 1).Value = $msgSheet.Cells.Item($r, 1).Value2
}

# Keep the originally active sheet (CypherOutput) selected, as it was before
# these sheets were appended.
$wb.Worksheets.Item("CypherOutput").Activate()
